$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("literature")

# ---- Row 12: Rahel, F. J. 2016 ... ----
# Copy formatting (style) from row 11 across A:D first, then overwrite values/formula.
$ws.Range("A11:D11").Copy() | Out-Null
$ws.Range("A12:D12").PasteSpecial(-4122) | Out-Null

$ws.Range("A12").Value = "2. Fisheries history"
$ws.Range("B12").Value = "Rahel, F. J. 2016. Changing Philosophies of Fisheries Management as Illustrated by the History of Fishing Regulations in Wyoming. Fisheries 41:38-48."
$ws.Range("C12").Value = "R189.pdf"

# ---- Row 13: Whelan, G. 2004 ... ----
# Copy formatting only to A13, C13, D13 (B13 is left at default/no explicit style).
$ws.Range("A11").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$ws.Range("C11").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$ws.Range("D11").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null

$ws.Range("A13").Value = "2. Fisheries history"
$ws.Range("B13").Value = "Whelan, G. 2004. A historical perspective on the use of propogated fish in fisheries management: Michigan's 130-year experience. American Fisheries Society Symposium 44:307-315."
$ws.Range("C13").Value = "W205.pdf"

# Fill the CONCATENATE formula down into the two new rows as one shared formula block.
$ws.Range("D12:D13").Formula = '=CONCATENATE(B12," [pdf](pdfs/",C12,")")'

# Row heights recalculated (wrapped-text autosize) for the edited/added rows.
$ws.Rows.Item(7).RowHeight = 72
$ws.Rows.Item(9).RowHeight = 57.6
$ws.Rows.Item(12).RowHeight = 72
$ws.Rows.Item(13).RowHeight = 57.6

# Update the view: scrolled down with A13 selected.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A13").Select() | Out-Null
